$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.446.67"
$ws.Range("E2").Value = "  +1.01%  "
$ws.Range("D3").Value = "3.140.64"
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "609.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.30%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "3.139.23"
$ws.Range("E8").Value = "  +0.34%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("E10").Value = "  +0.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.39"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.23%  "
$ws.Range("E12").Value = "  -0.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000258"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.42"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.42%  "
$ws.Range("D15").Value = "3.655.84"
$ws.Range("E15").Value = "  +0.27%  "
$ws.Range("E16").Value = "  +2.60%  "
$ws.Range("D17").Value = "64.393.76"
$ws.Range("E17").Value = "  +0.82%  "
$ws.Range("D18").Value = "3.140.73"
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.86"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "477.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.84"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.93%  "
$ws.Range("E22").Value = "  +2.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.76"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.64"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.41"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.96%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").Value = "  -3.42%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.44"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.80%  "
$ws.Range("E29").Value = "  +7.83%  "
$ws.Range("E30").Value = "  +2.27%  "
$ws.Range("E31").Value = "  -5.33%  "
$ws.Range("E32").Value = "  -0.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.86"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.64"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.11"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.99%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.99"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.70%  "
$ws.Range("D37").Value = "0.0₃0763"
$ws.Range("E37").Value = "  +4.42%  "
$ws.Range("E38").Value = "  -1.72%  "
$ws.Range("E39").Value = "  +3.89%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "445.98"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.87%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0393"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.40%  "
$ws.Range("E42").Value = "  +2.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.26"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "2.888.85"
$ws.Range("E44").Value = "  +1.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.262"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.24"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.42%  "
$ws.Range("E47").Value = "  +3.83%  "
$ws.Range("B48").Value = "USDe"
$ws.Range("C48").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.999"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "26.26"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.15%  "
$ws.Range("E50").Value = "  -0.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "119.65"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.51%  "
